$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update department names in column B (rename "Bilgisayar" -> "Bilgisayar Mühendisliği",
# "Endüstri" -> "Endüstri Mühendisliği")
$ws.Range("B2").Value = "Bilgisayar Mühendisliği"
$ws.Range("B3").Value = "Bilgisayar Mühendisliği"
$ws.Range("B4").Value = "Bilgisayar Mühendisliği"
$ws.Range("B5").Value = "Bilgisayar Mühendisliği"
$ws.Range("B6").Value = "Endüstri Mühendisliği"
$ws.Range("B7").Value = "Endüstri Mühendisliği"
$ws.Range("B8").Value = "Endüstri Mühendisliği"
$ws.Range("B9").Value = "Endüstri Mühendisliği"

# Fix "Günlük Max Ders Saati" (daily max course hours) values to all be 40
$ws.Range("D3").Value = 40
$ws.Range("D5").Value = 40
$ws.Range("D6").Value = 40
$ws.Range("D7").Value = 40
$ws.Range("D8").Value = 40

# Update the active selection shown in the sheet view
$ws.Range("F16").Select()
